$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the style of the existing date cell (A3) onto the new date cell (A4)
# so it reuses the same date/time number format instead of creating a new one.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)  # xlPasteFormats

# Add a new row of data (row 4) mirroring the existing rows' layout.
$ws.Cells.Item(4, 1).Value = 42606.881157407406
$ws.Cells.Item(4, 2).Value = 16
$ws.Cells.Item(4, 3).Value = 70
$ws.Cells.Item(4, 4).Value = 29
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 100
$ws.Cells.Item(4, 7).Value = 3421
$ws.Cells.Item(4, 8).Value = 3260
$ws.Cells.Item(4, 9).Value = 384
$ws.Cells.Item(4, 10).Value = 45
$ws.Cells.Item(4, 11).Value = 19
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 13).Value = 2
$ws.Cells.Item(4, 14).Value = "Bag"
